$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 90832

$ws.Range("A3").Value = 112206846
$ws.Range("B3").Value = 90806
$ws.Range("E3").Value = 4361
$ws.Range("F3").Value = "Orange taggsvamp"
$ws.Range("G3").Value = "Hydnellum aurantiacum"
$ws.Range("H3").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("AC3").Value = $null

$ws.Range("A4").Value = 112206831
$ws.Range("B4").Value = 90816
$ws.Range("E4").Value = 788
$ws.Range("F4").Value = "Gul taggsvamp"
$ws.Range("G4").Value = "Hydnellum geogenium"
$ws.Range("H4").Value = "(Fr.) Banker"
$ws.Range("AC4").Value = "Örtrikt dråg i granskog"

$ws.Range("A5").Value = 112395260
$ws.Range("B5").Value = 90832
$ws.Range("E5").Value = 4368
$ws.Range("F5").Value = "Dofttaggsvamp"
$ws.Range("G5").Value = "Hydnellum suaveolens"
$ws.Range("H5").Value = "(Scop.:Fr.) P. Karst."
$ws.Range("Q5").Value = 333038
$ws.Range("R5").Value = 6626631

$ws.Range("B6").Value = 90806

$ws.Range("B7").Value = 93307

$ws.Range("A8").Value = 112395257
$ws.Range("B8").Value = 90816
$ws.Range("E8").Value = 788
$ws.Range("F8").Value = "Gul taggsvamp"
$ws.Range("G8").Value = "Hydnellum geogenium"
$ws.Range("H8").Value = "(Fr.) Banker"
$ws.Range("Q8").Value = 333022
$ws.Range("R8").Value = 6626625
